# Weekly Fruta/Hortaliza update: two new daily price records (2023-03-09)
# are inserted at the top of the data block (rows 20-21), pushing the
# existing records down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 20, shifting rows 20:46 down to 22:48.
$ws.Rows("20:21").Insert()

$newDate = Get-Date -Year 2023 -Month 3 -Day 9 -Hour 0 -Minute 0 -Second 0

# Row 20: Membrillo Champion, calidad "Especial"
$ws.Cells.Item(20, 1).Value = 9
$ws.Cells.Item(20, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(20, 3).Value = "Metropolitana"
$ws.Cells.Item(20, 4).Value = $newDate
$ws.Cells.Item(20, 5).Value = 13
$ws.Cells.Item(20, 6).Value = "Fruta"
$ws.Cells.Item(20, 7).Value = 100104
$ws.Cells.Item(20, 8).Value = "Frutos de pepita"
$ws.Cells.Item(20, 9).Value = 100104003
$ws.Cells.Item(20, 10).Value = "Membrillo"
$ws.Cells.Item(20, 11).Value = "Champion"
$ws.Cells.Item(20, 12).Value = "Especial"
$ws.Cells.Item(20, 13).Value = 290
$ws.Cells.Item(20, 14).Value = 14400
$ws.Cells.Item(20, 15).Value = 14400
$ws.Cells.Item(20, 16).Value = 14400
$ws.Cells.Item(20, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(20, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(20, 19).Value = 800
$ws.Cells.Item(20, 20).Value = 18

# Row 21: Membrillo Champion, calidad "Primera"
$ws.Cells.Item(21, 1).Value = 9
$ws.Cells.Item(21, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(21, 3).Value = "Metropolitana"
$ws.Cells.Item(21, 4).Value = $newDate
$ws.Cells.Item(21, 5).Value = 13
$ws.Cells.Item(21, 6).Value = "Fruta"
$ws.Cells.Item(21, 7).Value = 100104
$ws.Cells.Item(21, 8).Value = "Frutos de pepita"
$ws.Cells.Item(21, 9).Value = 100104003
$ws.Cells.Item(21, 10).Value = "Membrillo"
$ws.Cells.Item(21, 11).Value = "Champion"
$ws.Cells.Item(21, 12).Value = "Primera"
$ws.Cells.Item(21, 13).Value = 350
$ws.Cells.Item(21, 14).Value = 12600
$ws.Cells.Item(21, 15).Value = 12600
$ws.Cells.Item(21, 16).Value = 12600
$ws.Cells.Item(21, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(21, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(21, 19).Value = 700
$ws.Cells.Item(21, 20).Value = 18
